# Edit Tab_7a_Quellen.xlsx:
#  - widen column G
#  - fix "Stifteverband" -> "Stifterverband" typo in row 44 (B:E)
#  - update the German/English homepage URLs in row 44 (F/G)
#  - add hyperlinks on F44/G44 pointing at those URLs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G (7th column) to fit the longer URL text.
# (52.4 is the ColumnWidth input that this engine's pixel-grid rounding
# maps closest to the target stored width of 53.171875.)
$ws.Columns.Item(7).ColumnWidth = 52.4

# Fix the misspelled organization name across B44:E44.
$fixedName = "Stifterverband Wissenschaftsstatistik"
$ws.Range("B44").Value = $fixedName
$ws.Range("C44").Value = $fixedName
$ws.Range("D44").Value = $fixedName
$ws.Range("E44").Value = $fixedName

# Update the German and English homepage URLs.
$deUrl = "https://www.stifterverband.org/"
$enUrl = "https://www.stifterverband.org/english"
$ws.Range("F44").Value = $deUrl
$ws.Range("G44").Value = $enUrl

# Add hyperlinks for the new URLs, mirroring the existing F3/G3 hyperlinks.
# Hyperlinks.Add re-styles the cell with Excel's built-in "Hyperlink" style
# by default, so capture the existing direct formatting first and reapply
# it afterwards (round-tripping the whole .Style object loses the cell's
# explicit formatting in this engine, so set the individual properties).
$cellF = $ws.Range("F44")
$cellG = $ws.Range("G44")

$fFontName = $cellF.Font.Name
$fFontSize = $cellF.Font.Size
$fFontUnderline = $cellF.Font.Underline
$fFontColor = $cellF.Font.Color
$fWrapText = $cellF.WrapText
$fVAlign = $cellF.VerticalAlignment
$fHAlign = $cellF.HorizontalAlignment
$fBorderColor = $cellF.Borders.Color
$fBorderStyle = $cellF.Borders.LineStyle

$gFontName = $cellG.Font.Name
$gFontSize = $cellG.Font.Size
$gFontUnderline = $cellG.Font.Underline
$gFontColor = $cellG.Font.Color
$gWrapText = $cellG.WrapText
$gVAlign = $cellG.VerticalAlignment
$gHAlign = $cellG.HorizontalAlignment
$gBorderColor = $cellG.Borders.Color
$gBorderStyle = $cellG.Borders.LineStyle

$ws.Hyperlinks.Add($cellF, $deUrl)
$ws.Hyperlinks.Add($cellG, $enUrl)

$cellF.Font.Name = $fFontName
$cellF.Font.Size = $fFontSize
$cellF.Font.Underline = $fFontUnderline
$cellF.Font.Color = $fFontColor
$cellF.WrapText = $fWrapText
$cellF.VerticalAlignment = $fVAlign
$cellF.HorizontalAlignment = $fHAlign
$cellF.Borders.Color = $fBorderColor
$cellF.Borders.LineStyle = $fBorderStyle

$cellG.Font.Name = $gFontName
$cellG.Font.Size = $gFontSize
$cellG.Font.Underline = $gFontUnderline
$cellG.Font.Color = $gFontColor
$cellG.WrapText = $gWrapText
$cellG.VerticalAlignment = $gVAlign
$cellG.HorizontalAlignment = $gHAlign
$cellG.Borders.Color = $gBorderColor
$cellG.Borders.LineStyle = $gBorderStyle
